$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("L2").Value = 3.5
$ws.Range("AC2").Value = 6
$ws.Range("AF2").Value = 81
$ws.Range("AS2").Value = 351
$ws.Range("AX2").Value = 17

# Row 3
$ws.Range("L3").Value = 4.5
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 2.75
$ws.Range("U3").Value = 2
$ws.Range("V3").Value = 1.73
$ws.Range("X3").Value = 8.5
$ws.Range("AG3").Value = 10
$ws.Range("AI3").Value = 15
$ws.Range("AO3").Value = 11

# Row 4
$ws.Range("G4").Value = 2.15
$ws.Range("H4").Value = 3.4
$ws.Range("I4").Value = 3.3
$ws.Range("O4").Value = 1.2
$ws.Range("P4").Value = 4.33
$ws.Range("Q4").Value = 1.7
$ws.Range("R4").Value = 2.1
$ws.Range("AH4").Value = 19
$ws.Range("AO4").Value = 11
$ws.Range("AP4").Value = 19
$ws.Range("AV4").Value = 41
$ws.Range("AW4").Value = 5.5

# Row 7
$ws.Range("G7").Value = 2.38
$ws.Range("H7").Value = 3.25
$ws.Range("I7").Value = 2.9
$ws.Range("J7").Value = 3.2
$ws.Range("L7").Value = 3.6
$ws.Range("N7").Value = 8.5
$ws.Range("X7").Value = 11
$ws.Range("Y7").Value = 10
$ws.Range("Z7").Value = 23
$ws.Range("AG7").Value = 8
$ws.Range("AH7").Value = 13
$ws.Range("AI7").Value = 11
$ws.Range("AJ7").Value = 29
$ws.Range("AK7").Value = 26
$ws.Range("AO7").Value = 15
$ws.Range("AQ7").Value = 51
$ws.Range("AW7").Value = 4.75
$ws.Range("AX7").Value = 17
$ws.Range("AZ7").Value = 51
$ws.Range("BB7").Value = 201
